$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-27 04:48:16'
$ws.Range('N2').Value = '0.5 °C 4:18 TU'
$ws.Range('E3').Value = '2026-02-27 04:48:18'
$ws.Range('N3').Value = '1.8 °C 4:23 TU'
$ws.Range('E4').Value = '2026-02-27 04:48:21'
$ws.Range('J4').Value = '1026.0 hPa'
$ws.Range('N4').Value = '6.3 °C 4:25 TU'
$ws.Range('O4').Value = '7.0 °C'
$ws.Range('E5').Value = '2026-02-27 04:48:23'
$ws.Range('H5').Value = '''38%'
$ws.Range('E6').Value = '2026-02-27 04:48:25'
$ws.Range('J6').Value = '1025.8 hPa'
$ws.Range('N6').Value = '8.9 °C 4:29 TU'
$ws.Range('E7').Value = '2026-02-27 04:48:28'
$ws.Range('J7').Value = '1026.1 hPa'
$ws.Range('N7').Value = '9.1 °C 4:23 TU'
$ws.Range('O7').Value = '10.2 °C'
$ws.Range('E8').Value = '2026-02-27 04:48:30'
$ws.Range('H8').Value = '''43%'
$ws.Range('J8').Value = '1025.5 hPa'
$ws.Range('L8').Value = '25.6 km/h - 258º 4:08 TU'
$ws.Range('O8').Value = '11.9 °C'
$ws.Range('E9').Value = '2026-02-27 04:48:33'
$ws.Range('M9').Value = '8.9 °C 4:28 TU'
$ws.Range('O9').Value = '8.3 °C'
$ws.Range('E10').Value = '2026-02-27 04:48:35'
$ws.Range('E11').Value = '2026-02-27 04:48:37'
$ws.Range('N11').Value = '1.7 °C 4:00 TU'
$ws.Range('O11').Value = '2.6 °C'
$ws.Range('E12').Value = '2026-02-27 04:48:39'
$ws.Range('M12').Value = '9.5 °C 4:18 TU'
$ws.Range('O12').Value = '7.8 °C'
$ws.Range('E13').Value = '2026-02-27 04:48:41'
$ws.Range('N13').Value = '-2.9 °C 4:29 TU'
$ws.Range('O13').Value = '-1.1 °C'
$ws.Range('E14').Value = '2026-02-27 04:48:42'
$ws.Range('N14').Value = '7.5 °C 4:29 TU'
$ws.Range('O14').Value = '9.3 °C'
$ws.Range('E15').Value = '2026-02-27 04:48:43'
$ws.Range('M15').Value = '9.4 °C 4:26 TU'
$ws.Range('O15').Value = '8.2 °C'
$ws.Range('E16').Value = '2026-02-27 04:48:44'
$ws.Range('H16').Value = '''25%'
$ws.Range('L16').Value = '31.3 km/h - 246º 4:20 TU'
$ws.Range('M16').Value = '3.7 °C 4:29 TU'
$ws.Range('O16').Value = '2.6 °C'
$ws.Range('E17').Value = '2026-02-27 04:48:45'
$ws.Range('E18').Value = '2026-02-27 04:48:46'
$ws.Range('N18').Value = '9.1 °C 4:24 TU'
$ws.Range('O18').Value = '9.9 °C'
$ws.Range('E19').Value = '2026-02-27 04:48:49'
$ws.Range('N19').Value = '6.9 °C 4:15 TU'
$ws.Range('E20').Value = '2026-02-27 04:48:51'
$ws.Range('H20').Value = '''61%'
$ws.Range('O20').Value = '1.7 °C'
$ws.Range('E21').Value = '2026-02-27 04:48:53'
$ws.Range('N21').Value = '2.6 °C 4:22 TU'
$ws.Range('O21').Value = '4.0 °C'
$ws.Range('E22').Value = '2026-02-27 04:48:56'
$ws.Range('H22').Value = '''48%'
$ws.Range('E23').Value = '2026-02-27 04:48:58'
$ws.Range('N23').Value = '1.6 °C 4:19 TU'
$ws.Range('O23').Value = '2.8 °C'
$ws.Range('E24').Value = '2026-02-27 04:49:00'
$ws.Range('N24').Value = '1.8 °C 4:20 TU'
$ws.Range('O24').Value = '4.9 °C'
$ws.Range('E25').Value = '2026-02-27 04:49:03'
$ws.Range('K25').Value = '-0.1 MJ/m2'
$ws.Range('O25').Value = '4.7 °C'
$ws.Range('E26').Value = '2026-02-27 04:49:05'
$ws.Range('H26').Value = '''46%'
$ws.Range('J26').Value = '1024.9 hPa'
$ws.Range('K26').Value = '-0.1 MJ/m2'
$ws.Range('M26').Value = '8.2 °C 4:26 TU'
$ws.Range('E27').Value = '2026-02-27 04:49:07'
$ws.Range('E28').Value = '2026-02-27 04:49:10'
$ws.Range('L28').Value = '6.8 km/h - 302º 4:12 TU'
$ws.Range('N28').Value = '4.5 °C 4:19 TU'
$ws.Range('O28').Value = '5.8 °C'
$ws.Range('E29').Value = '2026-02-27 04:49:12'
$ws.Range('L29').Value = '10.1 km/h - 202º 3:32 TU'
$ws.Range('M29').Value = '10.7 °C 3:35 TU'
$ws.Range('O29').Value = '9.5 °C'
$ws.Range('E30').Value = '2026-02-27 04:49:14'
$ws.Range('J30').Value = '1025.7 hPa'
$ws.Range('N30').Value = '9.6 °C 4:29 TU'
$ws.Range('O30').Value = '10.0 °C'
$ws.Range('E31').Value = '2026-02-27 04:49:16'
$ws.Range('L31').Value = '13.7 km/h - 334º 4:12 TU'
$ws.Range('O31').Value = '9.2 °C'
$ws.Range('E32').Value = '2026-02-27 04:49:19'
$ws.Range('H32').Value = '''92%'
$ws.Range('N32').Value = '0.5 °C 4:07 TU'
$ws.Range('O32').Value = '1.4 °C'
$ws.Range('E33').Value = '2026-02-27 04:49:21'
$ws.Range('N33').Value = '1.1 °C 4:27 TU'
$ws.Range('O33').Value = '2.7 °C'
$ws.Range('E34').Value = '2026-02-27 04:49:23'
$ws.Range('H34').Value = '''47%'
$ws.Range('L34').Value = '18.4 km/h - 13º 4:09 TU'
$ws.Range('M34').Value = '3.9 °C 4:07 TU'
$ws.Range('O34').Value = '1.8 °C'
$ws.Range('E35').Value = '2026-02-27 04:49:26'
$ws.Range('J35').Value = '1025.5 hPa'
$ws.Range('K35').Value = '-0.1 MJ/m2'
$ws.Range('O35').Value = '10.2 °C'
$ws.Range('E36').Value = '2026-02-27 04:49:28'
$ws.Range('M36').Value = '10.2 °C 4:14 TU'
$ws.Range('O36').Value = '9.0 °C'
$ws.Range('E37').Value = '2026-02-27 04:49:30'
$ws.Range('J37').Value = '1028.7 hPa'
$ws.Range('L37').Value = '19.1 km/h - 249º 4:15 TU'
$ws.Range('E38').Value = '2026-02-27 04:49:32'
$ws.Range('L38').Value = '11.2 km/h - 270º 4:03 TU'
$ws.Range('N38').Value = '6.9 °C 4:27 TU'
$ws.Range('O38').Value = '7.8 °C'
$ws.Range('E39').Value = '2026-02-27 04:49:35'
$ws.Range('O39').Value = '5.3 °C'
$ws.Range('E40').Value = '2026-02-27 04:49:37'
$ws.Range('H40').Value = '''98%'
$ws.Range('J40').Value = '1030.0 hPa'
$ws.Range('N40').Value = '1.0 °C 4:29 TU'
$ws.Range('O40').Value = '2.2 °C'
$ws.Range('E41').Value = '2026-02-27 04:49:39'
$ws.Range('J41').Value = '1026.1 hPa'
$ws.Range('N41').Value = '7.0 °C 4:18 TU'
$ws.Range('O41').Value = '9.1 °C'
$ws.Range('E42').Value = '2026-02-27 04:49:42'
$ws.Range('O42').Value = '8.2 °C'
$ws.Range('E43').Value = '2026-02-27 04:49:44'
$ws.Range('N43').Value = '2.8 °C 4:08 TU'
$ws.Range('O43').Value = '4.3 °C'
$ws.Range('E44').Value = '2026-02-27 04:49:47'
$ws.Range('H44').Value = '''68%'
$ws.Range('E45').Value = '2026-02-27 04:49:49'
$ws.Range('J45').Value = '1026.6 hPa'
$ws.Range('N45').Value = '5.1 °C 4:00 TU'
$ws.Range('O45').Value = '6.9 °C'
$ws.Range('E46').Value = '2026-02-27 04:49:51'
$ws.Range('J46').Value = '1026.3 hPa'
$ws.Range('N46').Value = '4.3 °C 4:29 TU'
$ws.Range('O46').Value = '7.1 °C'
